$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Untitled 1"
$ws.Range("C1").Value = "Untitled 2"
$ws.Range("D1").Value = "Untitled 3"

$ws.Range("A2").Value = 0.000000
$ws.Range("B2").Value = 0.000007
$ws.Range("C2").Value = -0.999635
$ws.Range("D2").Value = 12.562553

$ws.Range("A2").Copy() | Out-Null
$ws.Range("B2:D2").PasteSpecial(-4122) | Out-Null
